$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of regression-test results: "Baseline 2010-18 C106"
$row = 28

$ws.Cells.Item($row, 1).Value = "CW3M"
$ws.Cells.Item($row, 2).Value = "Baseline 2010-18 C106"
$ws.Cells.Item($row, 3).Value = "2010-18"

$ws.Cells.Item($row, 4).Value = 1186.9521077777779
$ws.Cells.Item($row, 5).Value = 1901.5157334444443
$ws.Cells.Item($row, 6).Value = 0.97970299999999988
$ws.Cells.Item($row, 7).Value = 280.33542888888883
$ws.Cells.Item($row, 8).Value = 9.775355222222224
$ws.Cells.Item($row, 9).Value = 12.968491888888888
$ws.Cells.Item($row, 10).Value = 8.145128999999999
$ws.Cells.Item($row, 11).Value = 645.94098588888892
$ws.Cells.Item($row, 12).Value = 83.47062044444445
$ws.Cells.Item($row, 13).Value = 1465.1962754444444
$ws.Cells.Item($row, 14).Value = 1191.1222331111112
$ws.Cells.Item($row, 15).Value = 4695.8937716666669
$ws.Cells.Item($row, 16).Value = 27227.338324888889
$ws.Cells.Item($row, 17).Value = 1.3484236666666667
$ws.Cells.Item($row, 18).Value = 0.00039399999999999998
$ws.Cells.Item($row, 19).Value = "2010-18"

# Match styling (number formats) used by the rest of the data rows
$ws.Range("D28:N28").NumberFormat = "0.00"
$ws.Range("O28:P28").NumberFormat = "0"
$ws.Range("Q28").NumberFormat = "0.00"
$ws.Range("R28").NumberFormat = "0.000000"

# Column I is highlighted (yellow fill) in this report, same as other rows
$ws.Range("I28").NumberFormat = "0.00"
$ws.Range("I28").Interior.Color = 65535

$ws.Range("S28").Select()
